# Update sheet name and data to reflect new "through" date (2022-10-14 -> 2022-10-15)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-10-15"

# Update the header label in I1 that shows the "through" date
$ws.Range("I1").Value = "2022 (through 10-15)"

# Update the data values that changed
$ws.Range("I11").Value = 48
$ws.Range("I14").Value = 1326
